$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Design")

# --- Insert a new column A ("Performer"/"Package") in front of the
#     existing "Workflow" table, shifting the old A:G columns to B:H ---
$ws.Columns.Item(1).Insert()

# --- Header row (row 1): add the new "Package" header in A1, clear the
#     old "Heading 1" (bold/border) row-level formatting so the header
#     row becomes plain text like the rest of the workbook ---
$ws.Range("A1").Value = "Package"
$ws.Rows.Item(1).ClearFormats()
$ws.Rows.Item(1).AutoFit()

# --- Row 2: "Performer" example row pulled in from the Example sheet ---
$ws.Range("A2").Value = "Performer"
$ws.Cells.Item(2,2).Value = "System_Module"
$ws.Cells.Item(2,3).Value = "text"
$ws.Cells.Item(2,4).Value = "text"
$ws.Cells.Item(2,5).Value = "text"
$ws.Cells.Item(2,6).Value = "in_Argument:string; out_Argument:int; io_Argument: dictionary"
$ws.Cells.Item(2,7).Value = "text"
$ws.Cells.Item(2,8).Value = "text"
$ws.Rows.Item(2).ClearFormats()
$ws.Rows.Item(2).AutoFit()
$ws.Range("B2:H2").Style = "Normal"
$ws.Range("B2:H2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 45

# --- Row 3: "Library" example row ---
$ws.Range("A3").Value = "Library"
$ws.Cells.Item(3,2).Value = "App_Module2"
$ws.Cells.Item(3,3).Value = "text"
$ws.Cells.Item(3,4).Value = "text"
$ws.Cells.Item(3,5).Value = "text"
$ws.Cells.Item(3,6).Value = "in_Argument:string; out_Argument:int; io_Argument: dictionary"
$ws.Cells.Item(3,7).Value = "text"
$ws.Cells.Item(3,8).Value = "text"
$ws.Rows.Item(3).ClearFormats()
$ws.Rows.Item(3).AutoFit()
$ws.Range("B3:H3").Style = "Normal"
$ws.Range("B3:H3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 45

# --- Row 4 stays blank; the new H4 cell already inherits the plain
#     wrap-text style ("s=1") from the column-insert shift, matching
#     the rest of the blank row. ---

# --- Remove the now-unused "Heading 1" cell style definition ---
$wb.Styles.Item("Heading 1").Delete()

# --- View/selection housekeeping to match the refreshed template ---
$ws.Activate()
$ws.Range("C11").Select()
$excel.ActiveWindow.Zoom = 115
